$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("scenarios")

# --- Update scenario values on the "scenarios" sheet ---
$ws1.Range("C4").Value = 5000000
$ws1.Range("C7").Value = 5000000

# --- Remove frozen pane, update selection ---
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("C19").Select() | Out-Null

# --- Add the new "sketch worksheet" sheet after "scenarios" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "sketch worksheet"
$ws2.PageSetup.FitToPagesWide = 1
$ws2.PageSetup.FitToPagesTall = $false

# Column widths (same as scenarios sheet)
$ws2.Columns.Item(2).ColumnWidth = 14.83203125
$ws2.Columns.Item(3).ColumnWidth = 14.83203125
$ws2.Columns.Item(4).ColumnWidth = 18.83203125
$ws2.Columns.Item(5).ColumnWidth = 25
$ws2.Columns.Item(6).ColumnWidth = 24.83203125

# Header row
$ws2.Range("A1").Value = "tp1"
$ws2.Range("B1").Value = "tree_analysis"
$ws2.Range("C1").Value = "draws"
$ws2.Range("D1").Value = "tree_final_states"
$ws2.Range("E1").Value = "damage_peak_temp"
$ws2.Range("F1").Value = "damage_disaster_tail"

# Data rows
$ws2.Range("A2").Value = 10
$ws2.Range("B2").Value = 4
$ws2.Range("C2").Value = 500000
$ws2.Range("D2").Value = 30
$ws2.Range("E2").Value = 11
$ws2.Range("F2").Value = 18
$ws2.Range("G2").Value = 45

$ws2.Range("A3").Value = 10
$ws2.Range("B3").Value = 4
$ws2.Range("C3").Value = 1000000
$ws2.Range("D3").Value = 30
$ws2.Range("E3").Value = 11
$ws2.Range("F3").Value = 18
$ws2.Range("G3").Formula = "=G2*C3/C2"

$ws2.Range("A4").Value = 10
$ws2.Range("B4").Value = 4
$ws2.Range("C4").Value = 5000000
$ws2.Range("D4").Value = 30
$ws2.Range("E4").Value = 11
$ws2.Range("F4").Value = 18
$ws2.Range("G4").Formula = "=G3*C4/C3"

$ws2.Range("A5").Value = 15
$ws2.Range("B5").Value = 4
$ws2.Range("C5").Value = 500000
$ws2.Range("D5").Value = 30
$ws2.Range("E5").Value = 11
$ws2.Range("F5").Value = 18
$ws2.Range("G5").Formula = "=G2*1.5"

$ws2.Range("A6").Value = 15
$ws2.Range("B6").Value = 4
$ws2.Range("C6").Value = 1000000
$ws2.Range("D6").Value = 30
$ws2.Range("E6").Value = 11
$ws2.Range("F6").Value = 18
$ws2.Range("G6").Formula = "=G3*1.5"

$ws2.Range("A7").Value = 15
$ws2.Range("B7").Value = 4
$ws2.Range("C7").Value = 5000000
$ws2.Range("D7").Value = 30
$ws2.Range("E7").Value = 11
$ws2.Range("F7").Value = 18
$ws2.Range("G7").Formula = "=G4*1.5"

Write-Output "done"
